{"js": "// Lattice multiplication worksheet: regenerate the exercise numbers.\n//\n// The document is a single 5-row x 3-column table. Each cell holds one\n// paragraph with one run whose text is split across 5 lines (joined by\n// manual line breaks, <w:br/>):\n//   line 1: \"A x B\"        (the two factors)\n//   line 2: \"  b1    b2\"   (digits of B, space padded -- needs xml:space\n//                           preserve since it starts with spaces)\n//   line 3: \"  ----\"       (separator, unchanged, needs xml:space preserve)\n//   line 4: \"a1|    |\"     (first digit of A, lattice grid)\n//   line 5: \"a2|    |\"     (second digit of A, lattice grid)\n//\n// We rewrite every cell's contents in place (same 5x3 grid shape before\n// and after the edit) using insertOoxml so we keep full control over the\n// run's <w:rPr> (font size 32) and the xml:space=\"preserve\" attribute on\n// just the lines that need it -- exactly mirroring how Word itself only\n// marks a text run \"preserve\" when it has leading/trailing whitespace.\n\nfunction buildCellOoxml(lines) {\n  const body = lines\n    .map((text, i) => {\n      const br = i > 0 ? \"<w:br/>\" : \"\";\n      const preserve = text !== text.trim() ? ' xml:space=\"preserve\"' : \"\";\n      const escaped = text\n        .replace(/&/g, \"&amp;\")\n        .replace(/</g, \"&lt;\")\n        .replace(/>/g, \"&gt;\");\n      return `${br}<w:t${preserve}>${escaped}</w:t>`;\n    })\n    .join(\"\");\n\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    `<w:p><w:r><w:rPr><w:sz w:val=\"32\"/></w:rPr>${body}</w:r></w:p>` +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\n// New contents for every cell, in row-major order (5 rows x 3 cols).\nconst newValues = [\n  [\"62 x 40\", \"  4    0\", \"  ----\", \"6|    |\", \"2|    |\"],\n  [\"18 x 50\", \"  5    0\", \"  ----\", \"1|    |\", \"8|    |\"],\n  [\"64 x 51\", \"  5    1\", \"  ----\", \"6|    |\", \"4|    |\"],\n\n  [\"13 x 31\", \"  3    1\", \"  ----\", \"1|    |\", \"3|    |\"],\n  [\"17 x 52\", \"  5    2\", \"  ----\", \"1|    |\", \"7|    |\"],\n  [\"22 x 58\", \"  5    8\", \"  ----\", \"2|    |\", \"2|    |\"],\n\n  [\"47 x 32\", \"  3    2\", \"  ----\", \"4|    |\", \"7|    |\"],\n  [\"67 x 84\", \"  8    4\", \"  ----\", \"6|    |\", \"7|    |\"],\n  [\"38 x 91\", \"  9    1\", \"  ----\", \"3|    |\", \"8|    |\"],\n\n  [\"58 x 31\", \"  3    1\", \"  ----\", \"5|    |\", \"8|    |\"],\n  [\"94 x 90\", \"  9    0\", \"  ----\", \"9|    |\", \"4|    |\"],\n  [\"34 x 40\", \"  4    0\", \"  ----\", \"3|    |\", \"4|    |\"],\n\n  [\"53 x 76\", \"  7    6\", \"  ----\", \"5|    |\", \"3|    |\"],\n  [\"86 x 63\", \"  6    3\", \"  ----\", \"8|    |\", \"6|    |\"],\n  [\"16 x 48\", \"  4    8\", \"  ----\", \"1|    |\", \"6|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < 5; r++) {\n  for (let c = 0; c < 3; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange(\"Whole\");\n    range.insertOoxml(buildCellOoxml(newValues[idx]), Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice multiplication worksheet: regenerate the exercise numbers.\n#\n# The document is a single 5-row x 3-column table. Each cell holds one\n# paragraph with one run whose text is split across 5 lines (joined by\n# manual line breaks, <w:br/>):\n#   line 1: \"A x B\"        (the two factors)\n#   line 2: \"  b1    b2\"   (digits of B, space padded -- needs xml:space\n#                           preserve since it starts with spaces)\n#   line 3: \"  ----\"       (separator, unchanged, needs xml:space preserve)\n#   line 4: \"a1|    |\"     (first digit of A, lattice grid)\n#   line 5: \"a2|    |\"     (second digit of A, lattice grid)\n#\n# We rewrite every cell's contents in place (same 5x3 grid shape before\n# and after the edit) using Range.InsertXML so we keep full control over\n# the run's <w:rPr> (font size 32) and the xml:space=\"preserve\" attribute\n# on just the lines that need it -- exactly mirroring how Word itself\n# only marks a text run \"preserve\" when it has leading/trailing\n# whitespace.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Build-CellXml([string[]]$lines) {\n    $wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'\n    $runXml = \"\"\n    for ($i = 0; $i -lt $lines.Length; $i++) {\n        $text = $lines[$i]\n        if ($i -gt 0) {\n            $runXml += \"<w:br/>\"\n        }\n        $escaped = $text.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n        if ($text -ne $text.Trim()) {\n            $runXml += '<w:t xml:space=\"preserve\">' + $escaped + '</w:t>'\n        } else {\n            $runXml += '<w:t>' + $escaped + '</w:t>'\n        }\n    }\n    return '<w:p xmlns:w=\"' + $wNs + '\"><w:r><w:rPr><w:sz w:val=\"32\"/></w:rPr>' + $runXml + '</w:r></w:p>'\n}\n\n# New contents for every cell, in row-major order (5 rows x 3 cols).\n$newValues = @(\n    @(\"62 x 40\", \"  4    0\", \"  ----\", \"6|    |\", \"2|    |\"),\n    @(\"18 x 50\", \"  5    0\", \"  ----\", \"1|    |\", \"8|    |\"),\n    @(\"64 x 51\", \"  5    1\", \"  ----\", \"6|    |\", \"4|    |\"),\n\n    @(\"13 x 31\", \"  3    1\", \"  ----\", \"1|    |\", \"3|    |\"),\n    @(\"17 x 52\", \"  5    2\", \"  ----\", \"1|    |\", \"7|    |\"),\n    @(\"22 x 58\", \"  5    8\", \"  ----\", \"2|    |\", \"2|    |\"),\n\n    @(\"47 x 32\", \"  3    2\", \"  ----\", \"4|    |\", \"7|    |\"),\n    @(\"67 x 84\", \"  8    4\", \"  ----\", \"6|    |\", \"7|    |\"),\n    @(\"38 x 91\", \"  9    1\", \"  ----\", \"3|    |\", \"8|    |\"),\n\n    @(\"58 x 31\", \"  3    1\", \"  ----\", \"5|    |\", \"8|    |\"),\n    @(\"94 x 90\", \"  9    0\", \"  ----\", \"9|    |\", \"4|    |\"),\n    @(\"34 x 40\", \"  4    0\", \"  ----\", \"3|    |\", \"4|    |\"),\n\n    @(\"53 x 76\", \"  7    6\", \"  ----\", \"5|    |\", \"3|    |\"),\n    @(\"86 x 63\", \"  6    3\", \"  ----\", \"8|    |\", \"6|    |\"),\n    @(\"16 x 48\", \"  4    8\", \"  ----\", \"1|    |\", \"6|    |\")\n)\n\n$idx = 0\nfor ($r = 1; $r -le 5; $r++) {\n    for ($c = 1; $c -le 3; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $xml = Build-CellXml $newValues[$idx]\n        [void]$cell.Range.InsertXML($xml)\n        $idx++\n    }\n}\n"}
